# Increment the "particip" (E) and "taxa_sucesso" (F) columns in the ranking
# table from a 0-1 fraction to a 0-100 scale (multiply each value by 100),
# for the data rows 2-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 7; $row++) {
    foreach ($col in @("E", "F")) {
        $cell = $ws.Range("$col$row")
        $cell.Value2 = $cell.Value2 * 100
    }
}
